$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Canada Premier League")

$row = 92
$prev = $row - 1

# Copy formatting from the row above so the new row matches existing styling
# (bold/bordered id column, date-formatted date column)
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($prev, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 90
$ws.Cells.Item($row, 2).Value = 7803362
$ws.Cells.Item($row, 3).Value = "Canada Premier League"
$ws.Cells.Item($row, 4).Value = "Canada Premier League"
$ws.Cells.Item($row, 5).Value = 45396.83333333334
$ws.Cells.Item($row, 6).Value = "Vancouver FC"
$ws.Cells.Item($row, 7).Value = "Valour FC"
$ws.Cells.Item($row, 8).Value = 4
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = "H"
$ws.Cells.Item($row, 11).Value = 2.4
$ws.Cells.Item($row, 12).Value = 3.6
$ws.Cells.Item($row, 13).Value = 2.4
$ws.Cells.Item($row, 14).Value = 2.9
$ws.Cells.Item($row, 15).Value = 3.5
$ws.Cells.Item($row, 16).Value = 2.05
$ws.Cells.Item($row, 17).Value = 0.5
$ws.Cells.Item($row, 18).Value = 1.8
$ws.Cells.Item($row, 19).Value = 2
$ws.Cells.Item($row, 20).Value = 2.5
$ws.Cells.Item($row, 21).Value = 1.95
$ws.Cells.Item($row, 22).Value = 1.85
$ws.Cells.Item($row, 23).Value = 1.9
$ws.Cells.Item($row, 24).Value = -1
$ws.Cells.Item($row, 25).Value = -1
$ws.Cells.Item($row, 26).Value = 0.8
$ws.Cells.Item($row, 27).Value = -1
$ws.Cells.Item($row, 28).Value = 0.95
$ws.Cells.Item($row, 29).Value = -1
